$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "2022-Q4" sheet right after "总计" ---------------------
# Copy the existing "2022-Q3" sheet (same column layout/styles as every
# quarter sheet) and drop the copy right after "总计"; this keeps the header
# row / borders / fonts identical to the other quarter tabs "for free".
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$sheetQ3.Copy($null, $sheetTotal)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Overwrite the copied data row with the new quarter's figures. The
# percentage-ish figures are stored as plain text (no special number
# format) in every quarter sheet, so write them with a leading quote to
# keep them text, then lift the formatting back off of an unstyled cell
# (H2) so no stray "quote prefix" style sticks around.
$newSheet.Range("D2").Value = "'0.21"
$newSheet.Range("E2").Value = "'73.13"
$newSheet.Range("F2").Value = "'3.07"
$newSheet.Range("G2").Value = "'0.0064"
$newSheet.Range("H2").Copy()
$newSheet.Range("D2:G2").PasteSpecial(-4122)
$newSheet.Range("H2").Value = 6

# --- 2. Update the "总计" summary sheet -------------------------------------
# Existing rows shift down one quarter (labels only - counts/values repeat),
# and a new row is appended for the quarter that has now rolled into view.
$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("B3").Value = "2022-Q3"
$sheetTotal.Range("B4").Value = "2022-Q2"

$sheetTotal.Range("A4").Copy()
$sheetTotal.Range("A5").PasteSpecial(-4122)
$sheetTotal.Range("A5").Value = 3
$sheetTotal.Range("B5").Value = "2022-Q1"
$sheetTotal.Range("C5").Value = 1
$sheetTotal.Range("D5").Value = 0.01

# Restore "2022-Q1" as the selected/active tab, matching the source state.
$wb.Worksheets.Item("2022-Q1").Activate()
